$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, styled like the other header cells (copy from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New "Save" column values H2:H9
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
